$wb = $excel.ActiveWorkbook

# --- Sheet "basic_search": add new column H with expected placeholder text ---
$ws1 = $wb.Worksheets.Item("basic_search")
$ws1.Range("H1").Value = "ExpectedPlaceholderText"
$ws1.Range("H1").Font.Bold = $true
$ws1.Range("H1").Font.Size = 11
$ws1.Range("H2:H7").Value = "Start typing to select a cancer type or keyword"
$ws1.Range("B7").Select()

# --- Sheet "basic_search_negative": add new columns F (CancerType) and G (Message) ---
$ws2 = $wb.Worksheets.Item("basic_search_negative")
$ws2.Range("F1").Value = "CancerType"
$ws2.Range("G1").Value = "Message"
$ws2.Range("F2").Value = "crab"
$ws2.Range("F3").Value = 123
$ws2.Range("F4").Value = "frye"
$ws2.Range("G2:G4").Value = "No available options found. Your search will be based on the text above."
$ws2.Columns.Item(7).ColumnWidth = 60.666666666666664

# --- Sheet "delighters": no content changes ---

# Final view state: "basic_search_negative" tab active with C12 selected
$ws2.Activate()
$ws2.Range("C12").Select()
